# C5-PowerPoint.pptx edit
# 1) Re-apply the (built-in) table style on the slide-6 table so it uses
#    {20786177-4114-420E-9C78-C0F02A797C9C} instead of the custom
#    {26FA4B71-32BA-4F87-8879-BF948F1A18E9} "Table_0" style.
# 2) Swap the deck's applied theme palette from "Integral" to the
#    standard Office "Office Theme" palette (dk1/lt1/dk2/lt2/accent1-6/
#    hlink/folHlink), matching the design-switch captured in the diff.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$tableSlide = $p.Slides.Item(6)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{20786177-4114-420E-9C78-C0F02A797C9C}")

# --- 2. Theme colours ------------------------------------------------------
# Order exposed by ThemeColorScheme: dk1, lt1, dk2, lt2, accent1..accent6,
# hlink, folHlink (12 entries) -- this writes straight into the theme part
# shared by every slide/layout/master in the deck.
$tcs = $p.Slides.Item(1).ThemeColorScheme
$tcs.Item(1).RGB  = 0         # dk1      = 000000
$tcs.Item(2).RGB  = 16777215  # lt1      = FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      = 44546A
$tcs.Item(4).RGB  = 15132391   # lt2      = E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  = 5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  = ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  = A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  = FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  = 4472C4
$tcs.Item(10).RGB = 4697456    # accent6  = 70AD47
$tcs.Item(11).RGB = 12673797   # hlink    = 0563C1
$tcs.Item(12).RGB = 7491477    # folHlink = 954F72
